# The deck's "Date Placeholder" (the datetimeFigureOut field shown in the
# lower-right of every slide) was re-cached by PowerPoint from 4/9/2023 to
# 12/30/2023. That cached text lives only on the Slide Master and on every
# Slide Layout (no individual slide overrides it), so update it there.

$p = $ppt.ActivePresentation
$newDate = "12/30/2023"

function Update-DatePlaceholders($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        # msoPlaceholder = 14, ppPlaceholderDate = 16
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

# Slide Master
Update-DatePlaceholders $p.SlideMaster.Shapes $newDate

# Every Slide Layout (CustomLayouts) hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}
